$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A189").Value = "TAO-USD"
$ws.Range("A190").Value = "IMX-USD"
$ws.Range("A191").Value = "GRT-USD"
